# Atualização das bases do grupo 13
#
# The sheet lists quarterly values of "Taxa de pessoas de 14 anos ou mais de
# idade, na força de trabalho, na semana de referência" for three regions,
# each stacked as a contiguous block of rows in chronological order:
#   Brasil   -> rows 2-24  (01/01/2019 .. 01/07/2024)
#   Nordeste -> rows 25-47 (01/01/2019 .. 01/07/2024)
#   Sergipe  -> rows 48-70 (01/01/2019 .. 01/07/2024)
#
# This update adds the newest quarter, 01/10/2024, to every block, which
# pushes the following blocks down by one row each time:
#   - new Brasil row   inserted at row 25 -> Nordeste/Sergipe shift by +1
#   - new Nordeste row inserted at row 49 -> Sergipe shifts by +1 more
#   - new Sergipe row  appended at the end, row 73

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variable = "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência"

# --- Brasil: insert new quarter 01/10/2024 right after the last Brasil row ---
$ws.Rows.Item(25).Insert()
$ws.Cells.Item(25, 1).Value2 = "Brasil"
$ws.Cells.Item(25, 2).Value2 = $variable
$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value2 = "01/10/2024"
$ws.Cells.Item(25, 4).Value2 = 93.83

# --- Nordeste: insert new quarter 01/10/2024 right after the last Nordeste row ---
$ws.Rows.Item(49).Insert()
$ws.Cells.Item(49, 1).Value2 = "Nordeste"
$ws.Cells.Item(49, 2).Value2 = $variable
$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 3).Value2 = "01/10/2024"
$ws.Cells.Item(49, 4).Value2 = 91.38

# --- Sergipe: append new quarter 01/10/2024 after the last Sergipe row ---
$ws.Cells.Item(73, 1).Value2 = "Sergipe"
$ws.Cells.Item(73, 2).Value2 = $variable
$ws.Cells.Item(73, 3).NumberFormat = "@"
$ws.Cells.Item(73, 3).Value2 = "01/10/2024"
$ws.Cells.Item(73, 4).Value2 = 91.47
